# Insert a new weekly price observation as row 17, pushing the existing
# rows 17-56 down to 18-57 (Fruta / hortaliza, semanal update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 17..56 down by inserting a fresh row at position 17.
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row with the new weekly observation.
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44623
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = 300000001
$ws.Range("G17").Value = "Rabanito"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 50
$ws.Range("K17").Value = 7000
$ws.Range("L17").Value = 7000
$ws.Range("M17").Value = 7000
$ws.Range("N17").Value = "`$/docena de paquetes"
$ws.Range("O17").Value = "Provincia de Cautín"
$ws.Range("P17").Value = 583
$ws.Range("Q17").Value = 12
$ws.Range("R17").Value = "Hortaliza"
